$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "Thanks for signing up..." text (A4) and append new rows.
# Values are assigned in the same order the shared-string table entries were
# originally authored in, so that the resulting xl/sharedStrings.xml ordering
# matches the target workbook exactly (the source rows were re-ordered after
# initial entry, which is why row order and shared-string order differ).

$ws.Range("A4").Value = "Search Results for: Virginia"
$ws.Range("A5").Value = "Housing Market"
$ws.Range("A6").Value = 139
$ws.Range("A7").Value = "What is the housing market like right now?"
$ws.Range("A8").Value = "22033 Apartments for Rent"
$ws.Range("A9").Value = "Dallas"
$ws.Range("A10").Value = "Miami"
$ws.Range("A11").Value = "Fairfax"
$ws.Range("A12").Value = "San Diego"
$ws.Range("A14").Value = "Orlando"
$ws.Range("A16").Value = "San Antonio"
$ws.Range("A17").Value = "Albuquerque"
$ws.Range("A18").Value = "New Orleans"
$ws.Range("A20").Value = "Wilmington"
$ws.Range("A15").Value = "San Francisco"
$ws.Range("A13").Value = "Las Vegas"
$ws.Range("A19").Value = "Boston"
$ws.Range("A21").Value = "Memphis"
$ws.Range("A22").Value = "Chicago"
$ws.Range("A23").Value = "Seattle"
$ws.Range("A24").Value = "Tampa"
$ws.Range("A25").Value = "Provo"
$ws.Range("A26").Value = "Denver"
$ws.Range("A27").Value = "Orem"
$ws.Range("A28").Value = "Phoenix"
$ws.Range("A29").Value = "Nashville"
$ws.Range("A30").Value = "Atlanta"
$ws.Range("A31").Value = "Boise"

# Apply the same formatting used on A3 ("Our Dedication to Data Quality") to
# the new "Housing Market" heading cell (paste formats only, so it reuses the
# existing cell style instead of defining a new one).
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A got narrower and the sheet view/selection moved to the new last
# cell, matching how Excel records the active cell after data entry.
$ws.Columns("A").ColumnWidth = 36.35
$null = $ws.Range("A31").Select()
